$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the report month in the title / subtitle text (October -> November)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Million Cubic Feet)"

# ---------------------------------------------------------------------------
# 2) Insert a new row for the "November" monthly figures (Year 2016 block).
#    This pushes the existing rows 53-60 down to 54-61 and keeps merged
#    ranges / row heights correctly shifted.
# ---------------------------------------------------------------------------
$ws.Rows.Item(53).Insert()

# Copy the formatting of the row above (October, row 52) down onto the new
# row so the label cell keeps style "8" and the data cells keep style "9",
# exactly like every other monthly row in this block.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new November row.
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 700215
$ws.Range("C53").Value = 332389
$ws.Range("D53").Value = 307969
$ws.Range("E53").Value = 5332
$ws.Range("F53").Value = 54526

# ---------------------------------------------------------------------------
# 3) Refresh the "Year to Date" annual block (now rows 54-57).
# ---------------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$ws.Range("B55").Value = 7870489
$ws.Range("C55").Value = 3591381
$ws.Range("D55").Value = 3645240
$ws.Range("E55").Value = 65957
$ws.Range("F55").Value = 567910

$ws.Range("B56").Value = 9209356
$ws.Range("C56").Value = 4351897
$ws.Range("D56").Value = 4225560
$ws.Range("E56").Value = 64344
$ws.Range("F56").Value = 567555

$ws.Range("B57").Value = 9698844
$ws.Range("C57").Value = 4702672
$ws.Range("D57").Value = 4344715
$ws.Range("E57").Value = 63156
$ws.Range("F57").Value = 588301

# ---------------------------------------------------------------------------
# 4) Refresh the "Rolling 12 Months Ending in November" block (now rows 59-60).
# ---------------------------------------------------------------------------
$ws.Range("B59").Value = 9883254
$ws.Range("C59").Value = 4655524
$ws.Range("D59").Value = 4534352
$ws.Range("E59").Value = 70344
$ws.Range("F59").Value = 623035

$ws.Range("B60").Value = 10506063
$ws.Range("C60").Value = 5096030
$ws.Range("D60").Value = 4695838
$ws.Range("E60").Value = 68904
$ws.Range("F60").Value = 645291
